# Keypad safe worst timing.xlsx - update:
#  - add a "combined" column (header E1 + formula E2 = D2/A2)
#  - A6's worst-case formula now multiplies by 288 instead of 280 and
#    uses the new "combined" value (E2) instead of D2
#  - A2's number format picks up the same font/alignment used by the
#    sheet's other header/body cells (style clean-up carried over from
#    the source workbook)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "combined" header cell (E1), formatted like the other headers (D1) ---
$ws.Range("E1").Value = "combined"
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)   # xlPasteFormats

# --- New E2 = D2 / A2, formatted like its neighbour D2 ---
$ws.Range("E2").Formula = "=D2/A2"
$ws.Range("D2").Copy()
$ws.Range("E2").PasteSpecial(-4122)   # xlPasteFormats

# --- A2 adopts the same format used elsewhere for labeled header cells (A1) ---
$ws.Range("A1").Copy()
$ws.Range("A2").PasteSpecial(-4122)   # xlPasteFormats

# --- Worst-case timing formula now derives from the combined value, 288 combos/sec ---
$ws.Range("A6").Formula = "=(E2*288) +20"

$excel.CutCopyMode = $false
